# Keyworddriven.xlsx - "Excel Sheet Final Changes"
#
# 1) Update the locator strings in column D (rows 4-9) to use the new
#    "name :" / "xpath :" separator style instead of "name=" / "xpath=".
# 2) Widen column D to fit the longer text.
# 3) Move the active selection from E7 to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "name :username"
$ws.Range("D5").Value = "name :password"
$ws.Range("D6").Value = 'xpath ://*[@id="app"]/div[1]/div/div[1]/div/div[2]/div[2]/form/div[3]/button'
$ws.Range("D7").Value = 'xpath ://*[@id=\"app\"]/div[1]/div[1]/aside/nav/div[2]/ul/li[9]/a'
$ws.Range("D8").Value = 'xpath ://*[@id=\"app\"]/div[1]/div[2]/div[2]/div/div[1]/div[2]/form/div[1]/div/div[2]/div/div[2]/div/div/div[1]'
$ws.Range("D9").Value = 'xpath ://*[@id=\"app\"]/div[1]/div[2]/div[2]/div/div[1]/div[2]/form/div[2]/button[2]'

$ws.Columns("D").ColumnWidth = 101.83333333333333

$ws.Range("F6").Select()
